$d = $word.ActiveDocument
Write-Output ("Content start=" + $d.Content.Start + " end=" + $d.Content.End)
$ftr1 = $d.Sections(1).Footers(1)
Write-Output ("Footer1 start=" + $ftr1.Range.Start + " end=" + $ftr1.Range.End)
$ftr2 = $d.Sections(1).Footers(2)
Write-Output ("Footer2 start=" + $ftr2.Range.Start + " end=" + $ftr2.Range.End)
$hdr1 = $d.Sections(1).Headers(1)
Write-Output ("Header1 exists=" + $hdr1.Exists + " start=" + $hdr1.Range.Start + " end=" + $hdr1.Range.End)
